# 02_02_Exercise.xlsx - apply the authored edits to the BEGIN sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BEGIN")

# --- Re-sort the "Data Point" / "Value" table (A2:B8) by Data Point ascending
# (the underlying Table5 sort condition moves from column B to column A).
# Styles stay fixed per row, so we only need to rewrite the values themselves
# in their final sorted order.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 20
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 200
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 20
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 140
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 40
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 20
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 50

# --- MEDIAN panel: clear the worked values that fed F2's formula
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()

# MEDIAN results row: F7 now references E4 (the "X Value" label) instead of F4
$ws.Range("F7").Formula = "=E4"

# --- MEAN panel: clear the worked values that fed J5's formula, and replace
# the duplicate "Mean" formula in J2 with a plain "Average" label
$ws.Range("J3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("J2").Value = "Average"

# --- MODE panel: clear the worked value, and repoint F9 at the I8 label
$ws.Range("J8").ClearContents()
$ws.Range("F9").Formula = "=I8"

# --- Update the active selection shown when the workbook is reopened
$ws.Range("G7").Select() | Out-Null
